# A new weekly price record was inserted at the top of the data table
# (row 91), and every subsequent record shifted down by one row, growing
# the table from 183 to 184 rows.
#
# Concretely:
#   - Row 91 keeps its identity but gets new observed values for the
#     "Fecha" (date, col D) and "Volumen" (col J) fields: the date moves
#     forward one day (44586 -> 44587) and the volume changes (560 -> 500).
#   - Rows 92..184 each take on the *previous* (pre-edit) row's values for
#     columns D,H,I,J,K,L,M,N,O,P,Q (Fecha, Variedad, Calidad, Volumen,
#     Precio minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg o
#     Unidades) - i.e. the whole table of price observations shifts down
#     by one, with a brand new row 184 appearing that holds what used to
#     be row 183's data.
#   - Columns A,B,C,E,F,G,R are constant across every data row (Mercado
#     ID, Mercado, Region, Codreg, Categoria ID, Categoria, Clasificacion)
#     so they are simply carried along unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row-to-row and that need to shift down.
$varyingCols = 4,8,9,10,11,12,13,14,15,16,17   # D,H,I,J,K,L,M,N,O,P,Q
# Columns that are constant across all data rows (copied along for safety).
$constCols   = 1,2,3,5,6,7,18                   # A,B,C,E,F,G,R

$firstDataRow = 91
$lastDataRow  = 183

# 1) Snapshot every current (pre-edit) value for rows 91..183 before any
#    writes happen, so the shifting writes below never read data we've
#    already overwritten.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowData = @{}
    foreach ($c in ($varyingCols + $constCols)) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $rowData['fmtD'] = $ws.Cells.Item($r, 4).NumberFormat
    $snapshot[$r] = $rowData
}

# 2) Shift rows 91..183 down into rows 92..184 (processing from the
#    bottom up so each destination row is written exactly once using the
#    snapshot captured above).
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + 1
    $rowData = $snapshot[$r]

    foreach ($c in $varyingCols) {
        $ws.Cells.Item($destRow, $c).Value = $rowData[$c]
    }
    foreach ($c in $constCols) {
        $ws.Cells.Item($destRow, $c).Value = $rowData[$c]
    }
    # Preserve the date display format on column D.
    $ws.Cells.Item($destRow, 4).NumberFormat = $rowData['fmtD']
}

# 3) Row 91 becomes the new observation: date advances one day and the
#    volume ("Volumen") value is updated; everything else it already had
#    (Variedad, Calidad, prices, unidad, origen, etc.) stays as-is.
$ws.Cells.Item($firstDataRow, 4).Value = 44587
$ws.Cells.Item($firstDataRow, 10).Value = 500
